$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19, 8).Value = 45455436
$ws.Cells.Item(19, 10).Value = 100000776
$ws.Cells.Item(19, 12).Value = 100000776
$ws.Cells.Item(19, 14).Value = -100001126
# Row 43
$ws.Cells.Item(43, 8).Value = 9793.134
$ws.Cells.Item(43, 9).Value = 9099
$ws.Cells.Item(43, 11).Value = 9099
$ws.Cells.Item(43, 13).Value = -9030
# Row 70
$ws.Cells.Item(70, 8).Value = 153190.58
$ws.Cells.Item(70, 10).Value = 13084
$ws.Cells.Item(70, 12).Value = 39252
$ws.Cells.Item(70, 14).Value = -39792
# Row 73
$ws.Cells.Item(73, 8).Value = 153190.58
$ws.Cells.Item(73, 10).Value = 13084
$ws.Cells.Item(73, 12).Value = 39252
$ws.Cells.Item(73, 14).Value = -41124
# Row 98
$ws.Cells.Item(98, 8).Value = 200754.03
$ws.Cells.Item(98, 10).Value = 1500377.5
$ws.Cells.Item(98, 12).Value = 1500377.5
$ws.Cells.Item(98, 14).Value = -1503373.5
# Row 122
$ws.Cells.Item(122, 8).Value = 200754.03
$ws.Cells.Item(122, 10).Value = 1500377.5
$ws.Cells.Item(122, 12).Value = 4501132.5
$ws.Cells.Item(122, 14).Value = -4506032.5
# Row 131
$ws.Cells.Item(131, 8).Value = 6805.2354
$ws.Cells.Item(131, 9).Value = 4616.3
$ws.Cells.Item(131, 10).Value = 9932.286
$ws.Cells.Item(131, 11).Value = 13848.9
$ws.Cells.Item(131, 12).Value = 29796.858
$ws.Cells.Item(131, 13).Value = -8808.900000000001
$ws.Cells.Item(131, 14).Value = -39876.858
# Row 138
$ws.Cells.Item(138, 8).Value = 2702.5059
$ws.Cells.Item(138, 9).Value = 1646.1428
$ws.Cells.Item(138, 10).Value = 3049.125
$ws.Cells.Item(138, 11).Value = 4938.428400000001
$ws.Cells.Item(138, 12).Value = 9147.375
$ws.Cells.Item(138, 13).Value = 201.5715999999993
$ws.Cells.Item(138, 14).Value = -19427.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 2807.9858
$ws.Cells.Item(32, 9).Value = 2154.7314
$ws.Cells.Item(32, 11).Value = 2154.7314
$ws.Cells.Item(32, 13).Value = -1867.7314
# Row 74
$ws.Cells.Item(74, 8).Value = 15155610
$ws.Cells.Item(74, 9).Value = 23811920
$ws.Cells.Item(74, 10).Value = 7067.5
$ws.Cells.Item(74, 11).Value = 23811920
$ws.Cells.Item(74, 12).Value = 7067.5
$ws.Cells.Item(74, 13).Value = -23811046
$ws.Cells.Item(74, 14).Value = -8815.5
# Row 77
$ws.Cells.Item(77, 8).Value = 15155610
$ws.Cells.Item(77, 9).Value = 23811920
$ws.Cells.Item(77, 10).Value = 7067.5
$ws.Cells.Item(77, 11).Value = 119059600
$ws.Cells.Item(77, 12).Value = 35337.5
$ws.Cells.Item(77, 13).Value = -119055232
$ws.Cells.Item(77, 14).Value = -44073.5
# Row 80
$ws.Cells.Item(80, 8).Value = 119997.5
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 11).Value = 0
$ws.Cells.Item(80, 13).ClearContents()
# Row 83
$ws.Cells.Item(83, 8).Value = 119997.5
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 11).Value = 0
$ws.Cells.Item(83, 13).ClearContents()
# Row 102
$ws.Cells.Item(102, 8).Value = 3165
$ws.Cells.Item(102, 9).Value = 1808.6
$ws.Cells.Item(102, 10).Value = 4295.3335
$ws.Cells.Item(102, 11).Value = 1808.6
$ws.Cells.Item(102, 12).Value = 4295.3335
$ws.Cells.Item(102, 13).Value = -186.5999999999999
$ws.Cells.Item(102, 14).Value = -7539.3335
# Row 110
$ws.Cells.Item(110, 8).Value = 2500.0688
$ws.Cells.Item(110, 9).Value = 1730.3462
$ws.Cells.Item(110, 11).Value = 1730.3462
$ws.Cells.Item(110, 13).Value = 314.6538
# Row 131
$ws.Cells.Item(131, 8).Value = 46065.168
$ws.Cells.Item(131, 10).Value = 46065.168
$ws.Cells.Item(131, 12).Value = 46065.168
$ws.Cells.Item(131, 14).Value = -56145.168
# Row 132
$ws.Cells.Item(132, 8).Value = 2183.8484
$ws.Cells.Item(132, 9).Value = 1550.0968
$ws.Cells.Item(132, 11).Value = 4650.2904
$ws.Cells.Item(132, 13).Value = -2120.2904

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 43367.035
$ws.Cells.Item(31, 9).Value = 5171.1055
$ws.Cells.Item(31, 11).Value = 5171.1055
$ws.Cells.Item(31, 13).Value = -4876.1055
# Row 34
$ws.Cells.Item(34, 8).Value = 43367.035
$ws.Cells.Item(34, 9).Value = 5171.1055
$ws.Cells.Item(34, 11).Value = 5171.1055
$ws.Cells.Item(34, 13).Value = -4969.1055
# Row 99
$ws.Cells.Item(99, 8).Value = 3738.4
$ws.Cells.Item(99, 9).Value = 3533.3333
$ws.Cells.Item(99, 11).Value = 3533.3333
$ws.Cells.Item(99, 13).Value = -2035.3333
# Row 126
$ws.Cells.Item(126, 8).Value = 3738.4
$ws.Cells.Item(126, 9).Value = 3533.3333
$ws.Cells.Item(126, 11).Value = 10599.9999
$ws.Cells.Item(126, 13).Value = -8129.999899999999
# Row 130
$ws.Cells.Item(130, 8).Value = 42361.875
$ws.Cells.Item(130, 10).Value = 50000
$ws.Cells.Item(130, 12).Value = 50000
$ws.Cells.Item(130, 14).Value = -60040
# Row 134
$ws.Cells.Item(134, 8).Value = 2702.1562
$ws.Cells.Item(134, 9).Value = 2223.3076
$ws.Cells.Item(134, 11).Value = 6669.9228
$ws.Cells.Item(134, 13).Value = -4134.9228

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 20
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 13).ClearContents()
# Row 42
$ws.Cells.Item(42, 8).Value = 9997.714
$ws.Cells.Item(42, 9).Value = 980
$ws.Cells.Item(42, 10).Value = 11500.667
$ws.Cells.Item(42, 11).Value = 2940
$ws.Cells.Item(42, 12).Value = 34502.001
$ws.Cells.Item(42, 13).Value = -2406
$ws.Cells.Item(42, 14).Value = -35570.001
# Row 130
$ws.Cells.Item(130, 8).Value = 1174
$ws.Cells.Item(130, 9).Value = 1174
$ws.Cells.Item(130, 10).Value = 0
$ws.Cells.Item(130, 11).Value = 3522
$ws.Cells.Item(130, 12).Value = 0
$ws.Cells.Item(130, 13).Value = 1498
$ws.Cells.Item(130, 14).ClearContents()
# Row 131
$ws.Cells.Item(131, 8).Value = 9682007
$ws.Cells.Item(131, 9).Value = 125000000
$ws.Cells.Item(131, 10).Value = 6078319.5
$ws.Cells.Item(131, 11).Value = 375000000
$ws.Cells.Item(131, 12).Value = 18234958.5
$ws.Cells.Item(131, 13).Value = -374994960
$ws.Cells.Item(131, 14).Value = -18245038.5
# Row 137
$ws.Cells.Item(137, 8).Value = 57122.055
$ws.Cells.Item(137, 9).Value = 834.44446
$ws.Cells.Item(137, 10).Value = 113409.664
$ws.Cells.Item(137, 11).Value = 2503.33338
$ws.Cells.Item(137, 12).Value = 340228.992
$ws.Cells.Item(137, 13).Value = 2596.66662
$ws.Cells.Item(137, 14).Value = -350428.992
# Row 141
$ws.Cells.Item(141, 8).Value = 6809.4375
$ws.Cells.Item(141, 9).Value = 3765.2727
$ws.Cells.Item(141, 11).Value = 11295.8181
$ws.Cells.Item(141, 13).Value = -6115.8181

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Cells.Item(21, 8).Value = 123399.4
$ws.Cells.Item(21, 9).Value = 352998.34
$ws.Cells.Item(21, 10).Value = 24999.857
$ws.Cells.Item(21, 11).Value = 352998.34
$ws.Cells.Item(21, 12).Value = 24999.857
$ws.Cells.Item(21, 13).Value = -352825.34
$ws.Cells.Item(21, 14).Value = -25345.857
# Row 30
$ws.Cells.Item(30, 8).Value = 123399.4
$ws.Cells.Item(30, 9).Value = 352998.34
$ws.Cells.Item(30, 10).Value = 24999.857
$ws.Cells.Item(30, 11).Value = 352998.34
$ws.Cells.Item(30, 12).Value = 24999.857
$ws.Cells.Item(30, 13).Value = -352893.34
$ws.Cells.Item(30, 14).Value = -25209.857
# Row 44
$ws.Cells.Item(44, 8).Value = 19998
$ws.Cells.Item(44, 10).Value = 19998
$ws.Cells.Item(44, 12).Value = 19998
$ws.Cells.Item(44, 14).Value = -21190
# Row 102
$ws.Cells.Item(102, 8).Value = 2083.8914
$ws.Cells.Item(102, 9).Value = 1568.1714
$ws.Cells.Item(102, 11).Value = 1568.1714
$ws.Cells.Item(102, 13).Value = 53.82860000000005
# Row 126
$ws.Cells.Item(126, 8).Value = 5854.727
$ws.Cells.Item(126, 9).Value = 4000
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 13).Value = -9530
# Row 132
$ws.Cells.Item(132, 8).Value = 3074.2942
$ws.Cells.Item(132, 9).Value = 1960.091
$ws.Cells.Item(132, 10).Value = 5117
$ws.Cells.Item(132, 11).Value = 5880.272999999999
$ws.Cells.Item(132, 12).Value = 15351
$ws.Cells.Item(132, 13).Value = -3350.272999999999
$ws.Cells.Item(132, 14).Value = -20411

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Cells.Item(40, 8).Value = 26518.5
$ws.Cells.Item(40, 9).Value = 32527.25
$ws.Cells.Item(40, 11).Value = 32527.25
$ws.Cells.Item(40, 13).Value = -32391.25
# Row 61
$ws.Cells.Item(61, 8).Value = 16977.8
$ws.Cells.Item(61, 9).Value = 20092.521
$ws.Cells.Item(61, 11).Value = 20092.521
$ws.Cells.Item(61, 13).Value = -19890.521
# Row 113
$ws.Cells.Item(113, 8).Value = 16977.8
$ws.Cells.Item(113, 9).Value = 20092.521
$ws.Cells.Item(113, 11).Value = 20092.521
$ws.Cells.Item(113, 13).Value = -17922.521

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 130
$ws.Cells.Item(130, 8).Value = 39462.75
$ws.Cells.Item(130, 10).Value = 39462.75
$ws.Cells.Item(130, 12).Value = 39462.75
$ws.Cells.Item(130, 14).Value = -49502.75
# Row 131
$ws.Cells.Item(131, 8).Value = 84062.5
$ws.Cells.Item(131, 10).Value = 84062.5
$ws.Cells.Item(131, 12).Value = 84062.5
$ws.Cells.Item(131, 14).Value = -94142.5
